# Insert a new data row at row 223 (pushing the existing rows 223-268 down
# to 224-269) and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 223:268 down by one row to make room for the new record.
$ws.Rows("223:223").Insert()

# Populate the newly inserted row 223 with the new observation.
$ws.Range("A223").Value = 4
$ws.Range("B223").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C223").Value = "Los Lagos"
$ws.Range("D223").Value = 44711
$ws.Range("E223").Value = 10
$ws.Range("F223").Value = 100112043
$ws.Range("G223").Value = "Pepino ensalada"
$ws.Range("H223").Value = "Sin especificar"
$ws.Range("I223").Value = "Primera"
$ws.Range("J223").Value = 120
$ws.Range("K223").Value = 25000
$ws.Range("L223").Value = 25000
$ws.Range("M223").Value = 25000
$ws.Range("N223").Value = "$/caja 60 unidades"
$ws.Range("O223").Value = "Región de Arica y Parinacota"
$ws.Range("P223").Value = 417
$ws.Range("Q223").Value = 60
$ws.Range("R223").Value = "Hortaliza"
